$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading History")

# Insert 5 new rows above the existing data (old row 5 -> new row 10, etc.)
$ws.Rows("5:9").Insert()

# Restore the date-style number format on the new rows' DATE column and
# clear the stray "header style" that Insert() propagated onto the rest
# of the newly created rows. (Kept to the same sparse columns Insert()
# itself materialized - A:N and W:AB - so we don't create any additional
# empty cells beyond what the row-shift already introduced.)
$ws.Range("A5:N9").ClearFormats()
$ws.Range("W5:AB9").ClearFormats()
$ws.Range("A5:A9").NumberFormat = "yyyy-mm-dd h:mm:ss"

# Row 5
$ws.Range("A5").Value = 45722
$ws.Range("B5").Value = "NSE"
$ws.Range("C5").Value = "Buy"
$ws.Range("D5").Value = 10
$ws.Range("E5").Value = 441.69
$ws.Range("F5").Value = 4416.9
$ws.Range("G5").Value = "~"
$ws.Range("J5").Formula = '=Index!$C$2'

# Row 6
$ws.Range("A6").Value = 45734
$ws.Range("B6").Value = "NSE"
$ws.Range("C6").Value = "Buy"
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = 469.7
$ws.Range("F6").Value = 2348.5
$ws.Range("G6").Value = "~"
$ws.Range("J6").Formula = '=Index!$C$2'

# Row 7
$ws.Range("A7").Value = 45811
$ws.Range("B7").Value = "NSE"
$ws.Range("C7").Value = "Buy"
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = 493.58
$ws.Range("F7").Value = 2467.9
$ws.Range("G7").Value = "~"
$ws.Range("J7").Formula = '=Index!$C$2'

# Row 8
$ws.Range("A8").Value = 45947
$ws.Range("B8").Value = "NSE"
$ws.Range("C8").Value = "Buy"
$ws.Range("D8").Value = 20
$ws.Range("E8").Value = 400.62
$ws.Range("F8").Value = 8012.4
$ws.Range("G8").Value = "~"
$ws.Range("J8").Formula = '=Index!$C$2'

# Row 9
$ws.Range("A9").Value = 45950
$ws.Range("B9").Value = "NSE"
$ws.Range("C9").Value = "Buy"
$ws.Range("D9").Value = 10
$ws.Range("E9").Value = 401.57
$ws.Range("F9").Value = 4015.7
$ws.Range("G9").Value = "~"
$ws.Range("J9").Formula = '=Index!$C$2'
